# Append a new day's log entry (row 44) to each of the four data sheets,
# mirroring the previous row's formatting/structure and advancing the
# timestamp in column A, per the latest scrape.

$wb = $excel.ActiveWorkbook

$newTimestamp = [double]"45830.43821759259"

foreach ($ws in $wb.Worksheets) {
    # Duplicate the last existing row (43) into the new row (44) so that
    # styles/number formats/cell types carry over exactly as they are.
    $ws.Range("A43:I43").Copy($ws.Range("A44:I44"))

    # The only column that actually changes day-over-day is the timestamp.
    $ws.Range("A44").Value = $newTimestamp
}
